$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.2588880590010376
$ws.Range("D2").Value = 592
$ws.Range("C3").Value = 0.5035661390002133
$ws.Range("D3").Value = 592
$ws.Range("C4").Value = 0.7889680520002003
$ws.Range("D4").Value = 582
$ws.Range("C5").Value = 0.9818535060003342
$ws.Range("D5").Value = 582
$ws.Range("C6").Value = 1.216101214999981
$ws.Range("D6").Value = 582
$ws.Range("C7").Value = 1.484270068000114
$ws.Range("D7").Value = 582
$ws.Range("C8").Value = 1.755511155000022
$ws.Range("D8").Value = 582
$ws.Range("C9").Value = 2.107468427999265
$ws.Range("D9").Value = 581
$ws.Range("C10").Value = 2.377706694998778
$ws.Range("D10").Value = 581
$ws.Range("C11").Value = 2.642980797998462
$ws.Range("D11").Value = 581
$ws.Range("C12").Value = 2.94413090299895
$ws.Range("D12").Value = 581
$ws.Range("C13").Value = 3.257306075999622
$ws.Range("D13").Value = 581
$ws.Range("C14").Value = 3.552105666999523
$ws.Range("D14").Value = 581
$ws.Range("C15").Value = 3.827662147999945
$ws.Range("D15").Value = 581
$ws.Range("C16").Value = 4.15952117300003
$ws.Range("D16").Value = 581
$ws.Range("C17").Value = 4.467262455000309
$ws.Range("D17").Value = 581
$ws.Range("C18").Value = 4.642935534000571
$ws.Range("D18").Value = 581
$ws.Range("C19").Value = 4.866894922000029
$ws.Range("D19").Value = 581
$ws.Range("C20").Value = 5.156737407000037
$ws.Range("D20").Value = 581
$ws.Range("C21").Value = 5.378081793999627
$ws.Range("D21").Value = 581
$ws.Range("C22").Value = 5.645050753999385
$ws.Range("D22").Value = 581
$ws.Range("C23").Value = 5.871087968998836
$ws.Range("D23").Value = 581
$ws.Range("C24").Value = 6.040054819998659
$ws.Range("D24").Value = 581
$ws.Range("C25").Value = 6.326487558998451
$ws.Range("D25").Value = 581
$ws.Range("C26").Value = 6.589295155998116
$ws.Range("D26").Value = 581
$ws.Range("C27").Value = 6.855491790998713
$ws.Range("D27").Value = 581
$ws.Range("C28").Value = 7.161866936998194
$ws.Range("D28").Value = 581
$ws.Range("C29").Value = 7.490124142997956
$ws.Range("D29").Value = 581
$ws.Range("C30").Value = 7.815030941997975
$ws.Range("D30").Value = 581
$ws.Range("C31").Value = 8.031460918998164
$ws.Range("D31").Value = 581
$ws.Range("C32").Value = 8.284992943998077
$ws.Range("D32").Value = 581
$ws.Range("C33").Value = 8.522040217998438
$ws.Range("D33").Value = 581
$ws.Range("C34").Value = 8.706479695998496
$ws.Range("D34").Value = 581
$ws.Range("C35").Value = 8.986557467997955
$ws.Range("D35").Value = 581
$ws.Range("C36").Value = 9.247064877998127
$ws.Range("D36").Value = 581
$ws.Range("C37").Value = 9.535354276998078
$ws.Range("D37").Value = 581
$ws.Range("C38").Value = 9.835632147997785
$ws.Range("D38").Value = 579
$ws.Range("C39").Value = 10.11193685399758
$ws.Range("D39").Value = 579
$ws.Range("C40").Value = 10.35034233999795
$ws.Range("D40").Value = 579
$ws.Range("C41").Value = 10.69565068599786
$ws.Range("D41").Value = 579
$ws.Range("C42").Value = 11.07080998499805
$ws.Range("D42").Value = 579
$ws.Range("C43").Value = 11.34120469099798
$ws.Range("D43").Value = 579
$ws.Range("C44").Value = 11.66284124999856
$ws.Range("D44").Value = 579
$ws.Range("C45").Value = 11.97010518399838
$ws.Range("D45").Value = 579
$ws.Range("C46").Value = 12.21920836599838
$ws.Range("D46").Value = 579
$ws.Range("C47").Value = 12.50621331599905
$ws.Range("D47").Value = 579
$ws.Range("C48").Value = 12.88723637799922
$ws.Range("D48").Value = 579
$ws.Range("C49").Value = 13.09753182599889
$ws.Range("D49").Value = 579
$ws.Range("C50").Value = 13.34658312499869
$ws.Range("D50").Value = 579
$ws.Range("C51").Value = 13.5534428539986
$ws.Range("D51").Value = 579
$ws.Range("C52").Value = 13.78137822199824
$ws.Range("D52").Value = 579
$ws.Range("C53").Value = 14.04427560099793
$ws.Range("D53").Value = 579
$ws.Range("C54").Value = 14.28526698099813
$ws.Range("D54").Value = 579
$ws.Range("C55").Value = 14.46571142599805
$ws.Range("D55").Value = 579
$ws.Range("C56").Value = 14.68069680399731
$ws.Range("D56").Value = 579
$ws.Range("C57").Value = 14.90298760299811
$ws.Range("D57").Value = 579
$ws.Range("C58").Value = 15.10509822299809
$ws.Range("D58").Value = 579
$ws.Range("C59").Value = 15.31814391799799
$ws.Range("D59").Value = 579
$ws.Range("C60").Value = 15.56520161399749
$ws.Range("D60").Value = 579
$ws.Range("C61").Value = 15.86397360399769
$ws.Range("D61").Value = 579
$ws.Range("C62").Value = 16.11630798899751
$ws.Range("D62").Value = 579
$ws.Range("C63").Value = 16.35316122399763
$ws.Range("D63").Value = 577
$ws.Range("C64").Value = 16.61892146399714
$ws.Range("D64").Value = 577
$ws.Range("C65").Value = 16.96522370799721
$ws.Range("D65").Value = 577
$ws.Range("C66").Value = 17.26415730199733
$ws.Range("D66").Value = 577
$ws.Range("C67").Value = 17.53253749299711
$ws.Range("D67").Value = 577
$ws.Range("C68").Value = 17.79850726099721
$ws.Range("D68").Value = 577
$ws.Range("C69").Value = 18.06145282899797
$ws.Range("D69").Value = 577
$ws.Range("C70").Value = 18.34682333199817
$ws.Range("D70").Value = 577
$ws.Range("C71").Value = 18.63151195299815
$ws.Range("D71").Value = 577
$ws.Range("C72").Value = 18.80927865399826
$ws.Range("D72").Value = 577
$ws.Range("C73").Value = 18.97913529699872
$ws.Range("D73").Value = 577
$ws.Range("C74").Value = 19.16493953999816
$ws.Range("D74").Value = 577
$ws.Range("C75").Value = 19.44286205099888
$ws.Range("D75").Value = 577
$ws.Range("C76").Value = 19.71336226899894
$ws.Range("D76").Value = 576
$ws.Range("C77").Value = 19.93776371899912
$ws.Range("D77").Value = 576
$ws.Range("C78").Value = 20.22123807099979
$ws.Range("D78").Value = 576
$ws.Range("C79").Value = 20.50536623599965
$ws.Range("D79").Value = 576
$ws.Range("C80").Value = 20.75574852199952
$ws.Range("D80").Value = 576
$ws.Range("C81").Value = 20.97906127799888
$ws.Range("D81").Value = 576
$ws.Range("C82").Value = 21.27304938499947
$ws.Range("D82").Value = 576
$ws.Range("C83").Value = 21.5467110489999
$ws.Range("D83").Value = 576
$ws.Range("C84").Value = 21.81063745199936
$ws.Range("D84").Value = 576
$ws.Range("C85").Value = 22.12620315999993
$ws.Range("D85").Value = 576
$ws.Range("C86").Value = 22.43145007000021
$ws.Range("D86").Value = 576
$ws.Range("C87").Value = 22.69962252799996
$ws.Range("D87").Value = 576
$ws.Range("C88").Value = 23.02606948799985
$ws.Range("D88").Value = 576
$ws.Range("C89").Value = 23.33156153999971
$ws.Range("D89").Value = 576
$ws.Range("C90").Value = 23.6409381200001
$ws.Range("D90").Value = 576
$ws.Range("C91").Value = 24.01073550000001
$ws.Range("D91").Value = 576
$ws.Range("C92").Value = 24.37368441099989
$ws.Range("D92").Value = 576
$ws.Range("C93").Value = 24.59619894000025
$ws.Range("D93").Value = 576
$ws.Range("C94").Value = 24.93442203300037
$ws.Range("D94").Value = 574
$ws.Range("C95").Value = 25.13742411600015
$ws.Range("D95").Value = 574
$ws.Range("C96").Value = 25.39841622800031
$ws.Range("D96").Value = 574
$ws.Range("C97").Value = 25.679469918
$ws.Range("D97").Value = 574
$ws.Range("C98").Value = 25.99392916500074
$ws.Range("D98").Value = 574
$ws.Range("C99").Value = 26.28218764100075
$ws.Range("D99").Value = 574
$ws.Range("C100").Value = 26.53933278200111
$ws.Range("D100").Value = 574
$ws.Range("C101").Value = 26.87106381200101
$ws.Range("D101").Value = 574
